# 自动更新Excel文件
# For every data row on the active sheet:
#   - Column D = total days ("总天")
#   - Column E = remaining days ("剩余")
#   - Column F = start date in yyyyMMdd form ("开始时间")
# Each day the remaining-day counter ticks down by one. When it would
# reach 0 (i.e. the previous remaining value was 1) the cycle restarts:
# remaining resets to the total-day count and the start date rolls
# forward by that many days.
# Rows whose start date cannot be parsed (corrupted data) are left
# untouched, just like the original nightly update job would skip them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal

    # Validate/parse the start date first; skip the whole row update
    # if the stored date is malformed (defensive, mirrors the
    # behaviour of the original automation script).
    $startStr = [string][int]$fVal
    try {
        $startDate = [datetime]::ParseExact($startStr, "yyyyMMdd", $null)
    }
    catch {
        continue
    }

    if ($remaining -gt 1) {
        $eCell.Value2 = $remaining - 1
    }
    elseif ($remaining -eq 1) {
        $eCell.Value2 = $totalDays
        $newDate = $startDate.AddDays($totalDays)
        $fCell.Value2 = [int]$newDate.ToString("yyyyMMdd")
    }
}
